# Applies the "Mejorar modal detalle (tabla, moneda, TC) y normalizar moneda" update
# to Bitacora_tareas.xlsx:
#  - Log: 4 new task rows (18-21)
#  - Resumen: 3 new rows inserted at row 15 + 1 new row appended at the end
#  - Versiones: 1 new version row (4)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Log": append rows 18-21
# ---------------------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")

$logRows = @(
    @("27/02/2025", "16:00", "Campo moneda en tabla transacciones", "Agregar columna moneda (ARS/USD) a la tabla transacciones en Supabase para normalizar la moneda de registración. Migración en supabase_transacciones_moneda.sql. Dashboard prioriza moneda; si viene vacío, infiere desde medio_pago. Export a Excel incluye columna moneda.", "Diagnostico"),
    @("27/02/2025", "16:20", "Modal detalle: ancho y moneda registración", "Ensanchado del modal mensual de detalle. En el listado de transacciones se muestra el monto con su moneda de registración (US$ / `$) antes del monto; si difiere de la moneda seleccionada, se muestra la conversión a la moneda de vista (→) o (sin cot.) si falta tipo de cambio.", "Diagnostico"),
    @("27/02/2025", "16:30", "Modal detalle: transacciones en tabla", "En el modal mensual (By Categoría / By Cuenta), el detalle expandido de transacciones ahora se renderiza como una tabla con encabezados (Fecha, Tipo, Medio, Mon., Monto, moneda vista, Descripción, Origen) para una lectura y análisis más clara.", "Diagnostico"),
    @("27/02/2025", "16:40", "Modal detalle: columna TC", "En la tabla de detalle expandida del modal mensual se agrega columna TC (MEP/CCL/OFICIAL según selector). Se muestra el tipo de cambio aplicado por fecha cuando hay conversión entre moneda de registración y moneda de vista; si no aplica muestra — y si falta cotización muestra sin cot.", "Diagnostico")
)

$startRow = 18
for ($i = 0; $i -lt $logRows.Count; $i++) {
    $r = $startRow + $i
    $row = $logRows[$i]
    $log.Cells.Item($r, 1).Value = $row[0]
    $log.Cells.Item($r, 2).Value = $row[1]
    $log.Cells.Item($r, 3).Value = $row[2]
    $log.Cells.Item($r, 4).Value = $row[3]
    $log.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------------
# Sheet "Resumen": insert 3 rows before row 15 (new "Detalle transacciones ..."
# entries), then append a new row at the end ("Campo moneda (BD)")
# ---------------------------------------------------------------------------
$resumen = $wb.Worksheets.Item("Resumen")

# Shift existing rows 15-24 down to 18-27, leaving rows 15-17 free.
$resumen.Rows.Item(15).Resize(3).Insert()

$newMidRows = @(
    @("Detalle transacciones (moneda registración)", "En el modal mensual, cada línea muestra el monto en su moneda de registración (US$ / `$). Si la moneda seleccionada difiere, se muestra la conversión a la moneda de vista (→) o indica (sin cot.) si falta tipo de cambio."),
    @("Detalle transacciones (tabla)", "En el modal mensual, al expandir una categoría/cuenta se muestra una tabla con títulos y filas de transacciones (Fecha, Tipo, Medio, Moneda, Monto, moneda vista, Descripción, Origen)."),
    @("Detalle transacciones (tipo de cambio)", "En el detalle expandido del modal mensual, se muestra la columna TC (según MEP/CCL/Oficial) cuando hay conversión entre moneda registración y moneda vista; si no aplica muestra — y si falta cotización muestra sin cot.")
)

for ($i = 0; $i -lt $newMidRows.Count; $i++) {
    $r = 15 + $i
    $row = $newMidRows[$i]
    $resumen.Cells.Item($r, 1).Value = $row[0]
    $resumen.Cells.Item($r, 2).Value = $row[1]
}

# New row appended at the end (row 28, after the former last row which is now 27)
$resumen.Cells.Item(28, 1).Value = "Campo moneda (BD)"
$resumen.Cells.Item(28, 2).Value = "Columna moneda en tabla transacciones (ARS/USD). Si está informada, el dashboard la usa; si no, infiere desde medio_pago (ej. ""dolar"" → USD). Export a Excel incluye moneda."

# ---------------------------------------------------------------------------
# Sheet "Versiones": append row 4 (version 1.2)
# ---------------------------------------------------------------------------
$versiones = $wb.Worksheets.Item("Versiones")

$versiones.Cells.Item(4, 1).Value = "1.2"
$versiones.Cells.Item(4, 2).Value = "27/02/2025"
$versiones.Cells.Item(4, 3).Value = "Modal mensual: detalle en tabla + moneda registración + TC; normalización moneda en BD y export Excel con moneda"
